$wb = $excel.ActiveWorkbook

# --- Work on the existing "Sheet4" sheet: clear the stray "-" marker in
# the Types 2 column for the first three data rows, and move the
# selection/active-cell away from B1 (this sheet is no longer the active
# tab once "cap" is added below). ---
$sheet4 = $wb.Worksheets.Item("Sheet4")
$sheet4.Activate()
$sheet4.Range("F2").Value = ""
$sheet4.Range("F3").Value = ""
$sheet4.Range("F4").Value = ""
$sheet4.Range("D8").Select()

# --- Add a new sheet "cap" right after "Sheet4" with a small market-cap
# snapshot (5 coins x 1 numeric row), and make it the active tab. ---
$capSheet = $wb.Worksheets.Add($null, $sheet4)
$capSheet.Name = "cap"

$capSheet.Range("A1").Value = "bitcoin"
$capSheet.Range("B1").Value = "ethereum"
$capSheet.Range("C1").Value = "ripple"
$capSheet.Range("D1").Value = "theta"
$capSheet.Range("E1").Value = "polkadot"

$capSheet.Range("A2:E2").WrapText = $true
$capSheet.Range("A2").Value = 870811968206
$capSheet.Range("B2").Value = 199308721408
$capSheet.Range("C2").Value = 20952203983
$capSheet.Range("D2").Value = 2897483033
$capSheet.Range("E2").Value = 20799036269

$capSheet.Range("A3").Select()
$capSheet.Activate()
